# Auto-generated script to update multiplication problems
$d = $word.ActiveDocument

$d.Content.Find.Execute("407×4=", $true, $false, $false, $false, $false, $true, 1, $false, "870×5=", 2) | Out-Null
$d.Content.Find.Execute("921×7=", $true, $false, $false, $false, $false, $true, 1, $false, "563×9=", 2) | Out-Null
$d.Content.Find.Execute("997×7=", $true, $false, $false, $false, $false, $true, 1, $false, "945×2=", 2) | Out-Null
$d.Content.Find.Execute("253×7=", $true, $false, $false, $false, $false, $true, 1, $false, "378×4=", 2) | Out-Null
$d.Content.Find.Execute("453×2=", $true, $false, $false, $false, $false, $true, 1, $false, "304×7=", 2) | Out-Null
$d.Content.Find.Execute("219×9=", $true, $false, $false, $false, $false, $true, 1, $false, "654×5=", 2) | Out-Null
$d.Content.Find.Execute("723×5=", $true, $false, $false, $false, $false, $true, 1, $false, "824×9=", 2) | Out-Null
$d.Content.Find.Execute("668×2=", $true, $false, $false, $false, $false, $true, 1, $false, "578×5=", 2) | Out-Null
$d.Content.Find.Execute("758×8=", $true, $false, $false, $false, $false, $true, 1, $false, "196×7=", 2) | Out-Null
$d.Content.Find.Execute("205×3=", $true, $false, $false, $false, $false, $true, 1, $false, "438×4=", 2) | Out-Null
$d.Content.Find.Execute("652×6=", $true, $false, $false, $false, $false, $true, 1, $false, "921×6=", 2) | Out-Null
$d.Content.Find.Execute("255×9=", $true, $false, $false, $false, $false, $true, 1, $false, "494×8=", 2) | Out-Null
$d.Content.Find.Execute("759×5=", $true, $false, $false, $false, $false, $true, 1, $false, "855×7=", 2) | Out-Null
$d.Content.Find.Execute("782×7=", $true, $false, $false, $false, $false, $true, 1, $false, "729×2=", 2) | Out-Null
$d.Content.Find.Execute("506×2=", $true, $false, $false, $false, $false, $true, 1, $false, "138×5=", 2) | Out-Null
$d.Content.Find.Execute("556×4=", $true, $false, $false, $false, $false, $true, 1, $false, "816×3=", 2) | Out-Null
$d.Content.Find.Execute("946×8=", $true, $false, $false, $false, $false, $true, 1, $false, "948×8=", 2) | Out-Null
$d.Content.Find.Execute("665×4=", $true, $false, $false, $false, $false, $true, 1, $false, "261×3=", 2) | Out-Null
$d.Content.Find.Execute("773×5=", $true, $false, $false, $false, $false, $true, 1, $false, "423×2=", 2) | Out-Null
$d.Content.Find.Execute("823×6=", $true, $false, $false, $false, $false, $true, 1, $false, "343×6=", 2) | Out-Null
$d.Content.Find.Execute("509×6=", $true, $false, $false, $false, $false, $true, 1, $false, "816×2=", 2) | Out-Null
$d.Content.Find.Execute("606×5=", $true, $false, $false, $false, $false, $true, 1, $false, "750×6=", 2) | Out-Null
$d.Content.Find.Execute("893×4=", $true, $false, $false, $false, $false, $true, 1, $false, "599×7=", 2) | Out-Null
$d.Content.Find.Execute("744×4=", $true, $false, $false, $false, $false, $true, 1, $false, "660×3=", 2) | Out-Null
$d.Content.Find.Execute("154×4=", $true, $false, $false, $false, $false, $true, 1, $false, "621×9=", 2) | Out-Null
